$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume figures (Price column D, Volume(1h) column E).
# NumberFormat is set to Text ("@") before assignment so that numeric-looking
# strings (e.g. "289.21") and percentages (e.g. "1.13%") are stored as literal
# text, matching the original inlineStr cell contents instead of being
# auto-converted to numbers/percentages by Excel.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "289.21"
$ws.Range("E2").Value = "1.13%"

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "29.40"
$ws.Range("E3").Value = "3.42%"

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.117"
$ws.Range("E4").Value = "3.77%"

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06678"
$ws.Range("E5").Value = "2.24%"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.80%"

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "3.411"
$ws.Range("E7").Value = "0.70%"

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "1.356"
$ws.Range("E8").Value = "-0.53%"

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9190"
$ws.Range("E9").Value = "0.94%"

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1587"
$ws.Range("E10").Value = "2.10%"

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06709"
$ws.Range("E11").Value = "-1.44%"

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07709"
$ws.Range("E12").Value = "0.66%"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.11%"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.36%"

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001568"
$ws.Range("E15").Value = "-2.03%"

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04514"
$ws.Range("E16").Value = "0.99%"

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0006446"
$ws.Range("E17").Value = "-1.30%"

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006270"
$ws.Range("E18").Value = "4.05%"

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "3.447"
$ws.Range("E19").Value = "-0.63%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.11%"

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3212"
$ws.Range("E21").Value = "1.98%"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.97%"

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "4.096"
$ws.Range("E23").Value = "1.19%"

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1569"
$ws.Range("E24").Value = "0.91%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.03%"

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004127"
$ws.Range("E26").Value = "-4.70%"

$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001251"
$ws.Range("E27").Value = "5.85%"

$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001619"
$ws.Range("E28").Value = "-1.05%"

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04215"
$ws.Range("E40").Value = "1.14%"

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006725"
$ws.Range("E41").Value = "0.43%"

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1240"
$ws.Range("E42").Value = "0.73%"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-8.48%"

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01217"
$ws.Range("E44").Value = "3.90%"

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005725"
$ws.Range("E45").Value = "3.11%"

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "1.972"
$ws.Range("E46").Value = "26.31%"

$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01308"
$ws.Range("E47").Value = "-29.34%"
